# Refatoração e ajuste de detalhes
# Merge the two product sheets ("produtoCategoria" + "produtoTexto") into a
# single sheet named "produtos" that carries Categoria, Produto and
# Pesquisa columns, updated with the new sample product/search data.

$wb = $excel.ActiveWorkbook

# The first sheet becomes the single surviving "produtos" sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "produtos"

# Add the third column (Pesquisa) header and refresh the sample row with
# the new product/category/search-term data.
$ws1.Range("A1").Value = "Produto"
$ws1.Range("B1").Value = "Categoria"
$ws1.Range("C1").Value = "Pesquisa"

$ws1.Range("A2").Value = "LOGITECH USB HEADSET H390"
$ws1.Range("B2").Value = "HEADPHONES"
$ws1.Range("C2").Value = "HEADSET"

# The second sheet ("produtoTexto") is no longer needed; remove it.
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("produtoTexto")
$ws2.Delete()
$excel.DisplayAlerts = $true

# Keep the same "next empty row" selection behaviour, now on column C.
$ws1.Activate()
$ws1.Range("C3").Select()
